# Refresh the crypto price/volume table (Price = column D, Volume(1h) = column E)
# with the latest scraped values. For column D we briefly force Text
# (NumberFormat "@") before assigning so numeric-looking strings such as
# "571.73", "0.110" or "0.0000178" are stored verbatim instead of being
# auto-parsed into a Number (which would lose trailing zeros or switch to
# scientific notation). ClearFormats() right after restores the cell to the
# workbook's default, unstyled appearance, matching the original cells
# (which carry no explicit style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.171.99"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.409.70"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("E4").Value = "  +0.65%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.73"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.35"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.58%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.538"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.431.49"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.110"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.37%  "
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.25"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.349"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.69"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000178"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +6.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.853.93"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.916.72"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.432.14"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.94"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.83"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.99"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("E22").Value = "  +2.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.03"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +12.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.30"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "614.97"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +10.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.42"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0988"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +7.97%  "
$ws.Range("E29").Value = "  +0.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.08"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.62%  "
$ws.Range("E31").Value = "  +8.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.83"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("E33").Value = "  +3.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.48"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.995"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.76"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "152.69"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("E38").Value = "  +1.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.39"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +5.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.53"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.68"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +16.58%  "
$ws.Range("E42").Value = "  +5.43%  "
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.18"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₆0282"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "143.43"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.58"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.32"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +6.85%  "
$ws.Range("E49").Value = "  +2.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0513"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0916"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.43%  "
